# Update the LR-pairs sheet with newly-recomputed TPM values.
# The new data set only has one row per "Sending cluster" (all targeting
# "Resolving-Mac"), so the old rows 6-9 (MuSCs / Resolving-Mac as sending
# cluster) are removed and rows 2-5 are rewritten in place with the
# updated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the last four data rows (old rows 6-9); the remaining rows 2-5
# shift up to become the new bottom of the table.
$ws.Range("A6:T9").EntireRow.Delete()

# New values for rows 2-5 (columns A..T).
$data = @(
    @("ECs","Icam4","Itgam","Resolving-Mac",3,1,1.829945333333333,5.489835999999999,0.4190796720210465,0.4190796720210465,3,1,35.68243999999999,107.04732,1,1,65.29691455994664,587.6722310395198,0.4190796720210465,0.4190796720210465),
    @("FAPs","Icam4","Itgam","Resolving-Mac",3,1,1.237199,3.711597,0.2833335737960661,0.2833335737960661,3,1,35.68243999999999,107.04732,1,1,44.14627908555999,397.3165117700399,0.2833335737960661,0.2833335737960661),
    @("MuSCs","Icam4","Itgam","Resolving-Mac",3,1,0.1530633333333333,0.45919,0.03505335944376924,0.03505335944376924,3,1,35.68243999999999,107.04732,1,1,5.461673207866665,49.15505887079999,0.03505335944376924,0.03505335944376924),
    @("Resolving-Mac","Icam4","Itgam","Resolving-Mac",3,1,1.146373333333333,3.43912,0.2625333947391181,0.2625333947391181,3,1,35.68243999999999,107.04732,1,1,40.90539768426665,368.1485791584,0.2625333947391181,0.2625333947391181)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowData[$j]
    }
}
